# Add a new row (29) of middle-school data to the sheet, matching the
# author's upload: six text-typed cells in columns A-F (G/H left blank).
#
# Columns A-C hold values that *look* numeric ("29", a latitude, a
# longitude). In the source workbook every data cell -- numeric-looking
# or not -- is stored as literal text (t="str"), so a plain
# `Range.Value = "29"` assignment (which Excel would auto-coerce to a
# real number) is wrong for those three. We force them to text by
# writing them as a formula that evaluates to a string, copying that
# result, and pasting-special as values only: the pasted result keeps
# the "text" type without requiring a `NumberFormat = "@"` style change
# (which would otherwise mint a new, unused cell style in styles.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue {
    param($range, [string]$text)

    # Scratch cell well outside the sheet's used range (A1:H28).
    $helper = $ws.Range("Z1")
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $excel.CutCopyMode = $false
    $helper.Clear()
}

$rowNum = 29

Set-TextCellValue $ws.Cells.Item($rowNum, 1) "29"
Set-TextCellValue $ws.Cells.Item($rowNum, 2) "34.34270611713106"
Set-TextCellValue $ws.Cells.Item($rowNum, 3) "133.9545587"

$ws.Cells.Item($rowNum, 4).Value = "高松市立下笠居中学校五色台分校"
$ws.Cells.Item($rowNum, 5).Value = "高松市中山町1501-192"
$ws.Cells.Item($rowNum, 6).Value = "087-811-6310"
